# Generate Report for Handback
# - "Status" of the ae87fc68... row changes from "Ready for handoff" to
#   "Handback transform failed" on every sheet that shows it.
# - The "Error Detail" column (P) for that same row now explains why the
#   handback transform failed, for both the zh-cn and de-de target sheets.
# - The Error Detail column is widened to fit the new message.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

$zhMessage = "Handback file name: btu4u1ll.of5 is different with handoff file name: ae87fc68-5935-4275-89f0-70c2f9351b2c.5ab37e874a09372f4bc638c58cf18b524f425c96.zh-cn."
$deMessage = "Handback file name: btu4u1ll.of5 is different with handoff file name: ae87fc68-5935-4275-89f0-70c2f9351b2c.5ab37e874a09372f4bc638c58cf18b524f425c96.de-de."

# --- Overview sheet: Status columns (E and F) for the ae87fc68 row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("P3").Value = $zhMessage
$wsZh.Range("P1").EntireColumn.ColumnWidth = 39.17

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("P3").Value = $deMessage
$wsDe.Range("P1").EntireColumn.ColumnWidth = 39.17
